# Reorder the "Requisitos" bullet list so that the Microbiologia
# (LOT2053) line moves from the first position to the last position,
# while Bioquimica I (LOT2007) and Engenharia Genetica (LOT2040) keep
# their relative order and each line stays its own run ending in a
# manual line break (<w:br/>).

$d = $word.ActiveDocument
$vt = [char]11   # vertical-tab == manual line break (<w:br/>) in Word's text model

$microbiologiaLine = "LOT2053 -  Microbiologia  (Requisito fraco)"

# Locate the "LOT2053 - Microbiologia ..." run, including its trailing
# line break, and remove it from the front of the list.
$matchRange = $d.Content
$found = $matchRange.Find.Execute($microbiologiaLine, $true, $false, $false,
                                   $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the Microbiologia requirement line"
}

$lineRange = $d.Range($matchRange.Start, $matchRange.End + 1)
$lineRange.Delete()

# Re-insert the same line (text + line break) as a brand-new run at the
# end of the "Requisitos" bullet paragraph (just before the paragraph
# mark), so it becomes the last item in the list.
$lastPara = $d.Paragraphs.Last
$endOfPara = $lastPara.Range.End
$insertPoint = $d.Range($endOfPara - 1, $endOfPara - 1)
$insertPoint.InsertAfter($microbiologiaLine + $vt)
